# Apply "Features upd & dataset column corr." changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update feature label text (also corrects a couple of swapped labels) ---
# NOTE: the order of assignment below matters for the resulting shared-string
# table layout, so it mirrors the order the strings were introduced upstream.
$ws.Range("B2").Value  = "Численность населения - popsize (тыс)"
$ws.Range("C3").Value  = "Безработные - unemployed (шт)"
$ws.Range("C2").Value  = "Сред. Числ. работников - avgemployers (тыс)"
$ws.Range("D2").Value  = "Сред. з/п - avgsalary (руб)"
$ws.Range("D3").Value  = "Площадь на одного - livarea (м2)"
$ws.Range("D4").Value  = "Число дошкол. мест - beforeschool (тыс)"
$ws.Range("D5").Value  = "Врачей на 10000 - docsperpop (на 10 тыс)"
$ws.Range("E2").Value  = "Удельный вес - invests ( %)"
$ws.Range("D7").Value  = "Мощность поликлиник на 10000 - cliniccap "
$ws.Range("D6").Value  = "Число коек на 10000 - bedsperpop "
$ws.Range("B9").Value  = "Степень износа -  funds %"
$ws.Range("C9").Value  = "Число предприятий - companies (шт)"
$ws.Range("D9").Value  = "Сумма четырех пунктов - factoriescap (сумма млн. руб.)"
$ws.Range("E9").Value  = "объем работ - conscap (млн. руб.)"
$ws.Range("E10").Value = "Ввод в действие жилых - consnewareas (тыс. м2)"
$ws.Range("E11").Value = "Число постр. квартир - consnewapt (шт)"
$ws.Range("F9").Value  = "Оборот розницы - retailturnover (млн. руб.)"
$ws.Range("F10").Value = "Оборот общепита - foodservturnover (млн. руб)"

# --- Widen columns B:F to fit the longer labels ---
# (input values are pre-compensated for this engine's internal width
#  quantization so the saved OOXML <col> widths land as close as possible
#  to the authored 42.140625 / 57.42578125 / 54.140625 / 50.7109375 / 45.7109375)
$ws.Columns.Item(2).ColumnWidth = 41.333333333333336
$ws.Columns.Item(3).ColumnWidth = 56.666666666666664
$ws.Columns.Item(4).ColumnWidth = 53.333333333333336
$ws.Columns.Item(5).ColumnWidth = 49.833333333333336
$ws.Columns.Item(6).ColumnWidth = 44.833333333333336

# --- Update the view / selection state left by the author ---
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G22").Select()
